$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The post "「恐れるな」" (row 225) was removed from the sheet.
# Deleting the entire row shifts all subsequent rows up by one,
# which matches the new dimension A1:C342.
$ws.Rows.Item(225).Delete()
